$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 125.083336
$ws.Range("I33").Value = 128.90909
$ws.Range("J33").Value = 83
$ws.Range("K33").Value = 128.90909
$ws.Range("L33").Value = 83
$ws.Range("M33").Value = 100.09091
$ws.Range("N33").Value = -541

$ws.Range("H106").Value = 2600.7144
$ws.Range("I106").Value = 2600.7144
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 2600.7144
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -1969.7144
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 143313.58
$ws.Range("I107").Value = 143313.58
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 143313.58
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -141393.58
$ws.Range("N107").ClearContents()

$ws.Range("H112").Value = 2572.5144
$ws.Range("I112").Value = 985
$ws.Range("J112").Value = 2668.7273
$ws.Range("K112").Value = 2955
$ws.Range("L112").Value = 8006.1819
$ws.Range("M112").Value = -1847
$ws.Range("N112").Value = -10222.1819

$ws.Range("H129").Value = 1069.381
$ws.Range("J129").Value = 1106.4177
$ws.Range("L129").Value = 3319.2531
$ws.Range("N129").Value = -13319.2531

$ws.Range("H137").Value = 2189.3696
$ws.Range("I137").Value = 1228.2858
$ws.Range("J137").Value = 2996.68
$ws.Range("K137").Value = 3684.8574
$ws.Range("L137").Value = 8990.039999999999
$ws.Range("M137").Value = -1134.8574
$ws.Range("N137").Value = -14090.04

$ws.Range("H138").Value = 2925.26
$ws.Range("J138").Value = 3313.0273
$ws.Range("L138").Value = 9939.081900000001
$ws.Range("N138").Value = -20219.0819

$ws.Range("H141").Value = 5132.4136
$ws.Range("I141").Value = 1771.8889
$ws.Range("J141").Value = 50499.5
$ws.Range("K141").Value = 5315.6667
$ws.Range("L141").Value = 151498.5
$ws.Range("M141").Value = -135.6666999999998
$ws.Range("N141").Value = -161858.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2740.5264
$ws.Range("I61").Value = 2751.3333
$ws.Range("J61").Value = 2700
$ws.Range("K61").Value = 2751.3333
$ws.Range("L61").Value = 2700
$ws.Range("M61").Value = -2539.3333
$ws.Range("N61").Value = -3124

$ws.Range("H74").Value = 1405.8776
$ws.Range("I74").Value = 1314.7368
$ws.Range("J74").Value = 1720.7273
$ws.Range("K74").Value = 1314.7368
$ws.Range("L74").Value = 1720.7273
$ws.Range("M74").Value = -440.7367999999999
$ws.Range("N74").Value = -3468.7273

$ws.Range("H77").Value = 1405.8776
$ws.Range("I77").Value = 1314.7368
$ws.Range("J77").Value = 1720.7273
$ws.Range("K77").Value = 6573.683999999999
$ws.Range("L77").Value = 8603.636500000001
$ws.Range("M77").Value = -2205.683999999999
$ws.Range("N77").Value = -17339.6365

$ws.Range("H132").Value = 4273.7886
$ws.Range("I132").Value = 4370.625
$ws.Range("J132").Value = 3951
$ws.Range("K132").Value = 13111.875
$ws.Range("L132").Value = 11853
$ws.Range("M132").Value = -10581.875
$ws.Range("N132").Value = -16913

$ws.Range("H136").Value = 2740.5264
$ws.Range("I136").Value = 2751.3333
$ws.Range("J136").Value = 2700
$ws.Range("K136").Value = 8253.999899999999
$ws.Range("L136").Value = 8100
$ws.Range("M136").Value = -5703.999899999999
$ws.Range("N136").Value = -13200

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 7145994
$ws.Range("I105").Value = 9526804
$ws.Range("J105").Value = 3564.2
$ws.Range("K105").Value = 9526804
$ws.Range("L105").Value = 3564.2
$ws.Range("M105").Value = -9525057
$ws.Range("N105").Value = -7058.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1769.381
$ws.Range("I31").Value = 2017.4651
$ws.Range("J31").Value = 1509.1951
$ws.Range("K31").Value = 2017.4651
$ws.Range("L31").Value = 1509.1951
$ws.Range("M31").Value = -1722.4651
$ws.Range("N31").Value = -2099.1951

$ws.Range("H34").Value = 1769.381
$ws.Range("I34").Value = 2017.4651
$ws.Range("J34").Value = 1509.1951
$ws.Range("K34").Value = 2017.4651
$ws.Range("L34").Value = 1509.1951
$ws.Range("M34").Value = -1815.4651
$ws.Range("N34").Value = -1913.1951

$ws.Range("H134").Value = 1538.4348
$ws.Range("I134").Value = 1260.4762
$ws.Range("J134").Value = 4457
$ws.Range("K134").Value = 3781.4286
$ws.Range("L134").Value = 13371
$ws.Range("M134").Value = -1246.4286
$ws.Range("N134").Value = -18441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 2000
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 6000
$ws.Range("N32").Value = -6566
$ws.Range("M32").ClearContents()

$ws.Range("H57").Value = 2789
$ws.Range("I57").Value = 2789
$ws.Range("K57").Value = 8367
$ws.Range("M57").Value = -7808

$ws.Range("H113").Value = 586.17645
$ws.Range("I113").Value = 586.17645
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1758.52935
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 411.4706499999998
$ws.Range("N113").ClearContents()

$ws.Range("H129").Value = 4167159
$ws.Range("I129").Value = 288.57144
$ws.Range("J129").Value = 10000778
$ws.Range("K129").Value = 865.71432
$ws.Range("L129").Value = 30002334
$ws.Range("M129").Value = 4134.28568
$ws.Range("N129").Value = -30012334

$ws.Range("H133").Value = 4684.375
$ws.Range("J133").Value = 6291.8
$ws.Range("L133").Value = 18875.4
$ws.Range("N133").Value = -28995.4

$ws.Range("H134").Value = 3373.5908
$ws.Range("I134").Value = 2351.9333
$ws.Range("J134").Value = 5562.857
$ws.Range("K134").Value = 7055.7999
$ws.Range("L134").Value = 16688.571
$ws.Range("M134").Value = -1985.7999
$ws.Range("N134").Value = -26828.571

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 32000
$ws.Range("J62").Value = 32000
$ws.Range("L62").Value = 32000
$ws.Range("N62").Value = -33372

$ws.Range("H65").Value = 32000
$ws.Range("J65").Value = 32000
$ws.Range("L65").Value = 96000
$ws.Range("N65").Value = -102864

$ws.Range("H122").Value = 1017
$ws.Range("I122").Value = 871.3333
$ws.Range("J122").Value = 1454
$ws.Range("K122").Value = 2613.9999
$ws.Range("L122").Value = 4362
$ws.Range("M122").Value = -163.9998999999998
$ws.Range("N122").Value = -9262

$ws.Range("H126").Value = 2316.7058
$ws.Range("I126").Value = 1376.25
$ws.Range("J126").Value = 3152.6667
$ws.Range("K126").Value = 4128.75
$ws.Range("L126").Value = 9458.000100000001
$ws.Range("M126").Value = -1658.75
$ws.Range("N126").Value = -14398.0001

$ws.Range("H132").Value = 3855.7896
$ws.Range("I132").Value = 3483
$ws.Range("J132").Value = 4899.6
$ws.Range("K132").Value = 10449
$ws.Range("L132").Value = 14698.8
$ws.Range("M132").Value = -7919
$ws.Range("N132").Value = -19758.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1264.8572
$ws.Range("I46").Value = 1077.5
$ws.Range("J46").Value = 1514.6666
$ws.Range("K46").Value = 1077.5
$ws.Range("L46").Value = 1514.6666
$ws.Range("M46").Value = -889.5
$ws.Range("N46").Value = -1890.6666

$ws.Range("H132").Value = 4689.357
$ws.Range("I132").Value = 4828.2334
$ws.Range("J132").Value = 4342.1665
$ws.Range("K132").Value = 14484.7002
$ws.Range("L132").Value = 13026.4995
$ws.Range("M132").Value = -11954.7002
$ws.Range("N132").Value = -18086.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 10000
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 10000
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

$ws.Range("H107").Value = 689.6667
$ws.Range("I107").Value = 376
$ws.Range("J107").Value = 940.6
$ws.Range("K107").Value = 1128
$ws.Range("L107").Value = 2821.8
$ws.Range("M107").Value = 792
$ws.Range("N107").Value = -6661.8

$ws.Range("H132").Value = 3519.2068
$ws.Range("I132").Value = 3737.3333
$ws.Range("J132").Value = 3285.5
$ws.Range("K132").Value = 11211.9999
$ws.Range("L132").Value = 9856.5
$ws.Range("M132").Value = -8681.999899999999
$ws.Range("N132").Value = -14916.5

$ws.Range("H141").Value = 77133.75
$ws.Range("J141").Value = 77133.75
$ws.Range("L141").Value = 77133.75
$ws.Range("N141").Value = -87493.75
